$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 436.66666
$ws.Cells.Item(9, 9).Value = 307.5
$ws.Cells.Item(9, 10).Value = 540
$ws.Cells.Item(9, 11).Value = 307.5
$ws.Cells.Item(9, 12).Value = 540
$ws.Cells.Item(9, 13).Value = -138.5
$ws.Cells.Item(9, 14).Value = -878

$ws.Cells.Item(17, 8).Value = 1952.5
$ws.Cells.Item(17, 10).Value = 1907.1111
$ws.Cells.Item(17, 12).Value = 5721.3333
$ws.Cells.Item(17, 14).Value = -6057.3333

$ws.Cells.Item(28, 8).Value = 981.3889
$ws.Cells.Item(28, 9).Value = 936.7692
$ws.Cells.Item(28, 11).Value = 936.7692
$ws.Cells.Item(28, 13).Value = -451.7692

$ws.Cells.Item(32, 8).Value = 16669967
$ws.Cells.Item(32, 9).Value = 1100
$ws.Cells.Item(32, 10).Value = 20003740
$ws.Cells.Item(32, 11).Value = 1100
$ws.Cells.Item(32, 12).Value = 20003740
$ws.Cells.Item(32, 13).Value = -774
$ws.Cells.Item(32, 14).Value = -20004392

$ws.Cells.Item(112, 8).Value = 1665.9714
$ws.Cells.Item(112, 9).Value = 1232.5
$ws.Cells.Item(112, 10).Value = 1692.2424
$ws.Cells.Item(112, 11).Value = 3697.5
$ws.Cells.Item(112, 12).Value = 5076.7272
$ws.Cells.Item(112, 13).Value = -2589.5
$ws.Cells.Item(112, 14).Value = -7292.7272

$ws.Cells.Item(113, 8).Value = 10272.223
$ws.Cells.Item(113, 9).Value = 7290
$ws.Cells.Item(113, 10).Value = 14000
$ws.Cells.Item(113, 11).Value = 7290
$ws.Cells.Item(113, 12).Value = 14000
$ws.Cells.Item(113, 13).Value = -4036
$ws.Cells.Item(113, 14).Value = -20508

$ws.Cells.Item(121, 8).Value = 3000
$ws.Cells.Item(121, 10).Value = 3000
$ws.Cells.Item(121, 12).Value = 9000
$ws.Cells.Item(121, 14).Value = -12494

$ws.Cells.Item(125, 8).Value = 35346.832
$ws.Cells.Item(125, 9).Value = 51745.25
$ws.Cells.Item(125, 11).Value = 465707.25
$ws.Cells.Item(125, 13).Value = -463247.25

$ws.Cells.Item(129, 8).Value = 252794.25
$ws.Cells.Item(129, 10).Value = 1836.6
$ws.Cells.Item(129, 12).Value = 5509.799999999999
$ws.Cells.Item(129, 14).Value = -15509.8

$ws.Cells.Item(131, 8).Value = 10139.087
$ws.Cells.Item(131, 9).Value = 6893.2666
$ws.Cells.Item(131, 11).Value = 20679.7998
$ws.Cells.Item(131, 13).Value = -15639.7998

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 8581.35
$ws.Cells.Item(32, 9).Value = 9287.638999999999
$ws.Cells.Item(32, 11).Value = 9287.638999999999
$ws.Cells.Item(32, 13).Value = -9000.638999999999

$ws.Cells.Item(74, 8).Value = 1474.326
$ws.Cells.Item(74, 9).Value = 1488.2059
$ws.Cells.Item(74, 11).Value = 1488.2059
$ws.Cells.Item(74, 13).Value = -614.2058999999999

$ws.Cells.Item(77, 8).Value = 1474.326
$ws.Cells.Item(77, 9).Value = 1488.2059
$ws.Cells.Item(77, 11).Value = 7441.0295
$ws.Cells.Item(77, 13).Value = -3073.0295

$ws.Cells.Item(88, 8).Value = 1888.6
$ws.Cells.Item(88, 9).Value = 2035.3334
$ws.Cells.Item(88, 10).Value = 1668.5
$ws.Cells.Item(88, 11).Value = 2035.3334
$ws.Cells.Item(88, 12).Value = 1668.5
$ws.Cells.Item(88, 13).Value = -1629.3334
$ws.Cells.Item(88, 14).Value = -2480.5

$ws.Cells.Item(91, 8).Value = 1888.6
$ws.Cells.Item(91, 9).Value = 2035.3334
$ws.Cells.Item(91, 10).Value = 1668.5
$ws.Cells.Item(91, 11).Value = 2035.3334
$ws.Cells.Item(91, 12).Value = 1668.5
$ws.Cells.Item(91, 13).Value = -631.3334
$ws.Cells.Item(91, 14).Value = -4476.5

$ws.Cells.Item(132, 8).Value = 1766.5714
$ws.Cells.Item(132, 9).Value = 1634.7778
$ws.Cells.Item(132, 11).Value = 4904.3334
$ws.Cells.Item(132, 13).Value = -2374.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 18520118
$ws.Cells.Item(86, 9).Value = 20001592
$ws.Cells.Item(86, 11).Value = 20001592
$ws.Cells.Item(86, 13).Value = -20000469

$ws.Cells.Item(89, 8).Value = 18520118
$ws.Cells.Item(89, 9).Value = 20001592
$ws.Cells.Item(89, 11).Value = 100007960
$ws.Cells.Item(89, 13).Value = -100002344

$ws.Cells.Item(94, 8).Value = 1932.5405
$ws.Cells.Item(94, 9).Value = 1696.56
$ws.Cells.Item(94, 11).Value = 1696.56
$ws.Cells.Item(94, 13).Value = -1245.56

$ws.Cells.Item(134, 8).Value = 3180.6667
$ws.Cells.Item(134, 9).Value = 2901.25
$ws.Cells.Item(134, 11).Value = 8703.75
$ws.Cells.Item(134, 13).Value = -6168.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 325.625
$ws.Cells.Item(7, 9).Value = 36.666668
$ws.Cells.Item(7, 10).Value = 499
$ws.Cells.Item(7, 11).Value = 36.666668
$ws.Cells.Item(7, 12).Value = 499
$ws.Cells.Item(7, 13).Value = 76.333332
$ws.Cells.Item(7, 14).Value = -725

$ws.Cells.Item(58, 8).Value = 2026.15
$ws.Cells.Item(58, 9).Value = 1166.6666
$ws.Cells.Item(58, 10).Value = 2729.3635
$ws.Cells.Item(58, 11).Value = 1166.6666
$ws.Cells.Item(58, 12).Value = 2729.3635
$ws.Cells.Item(58, 13).Value = -963.6666
$ws.Cells.Item(58, 14).Value = -3135.3635

$ws.Cells.Item(99, 8).Value = 2931.476
$ws.Cells.Item(99, 9).Value = 3018.9167
$ws.Cells.Item(99, 10).Value = 2814.889
$ws.Cells.Item(99, 11).Value = 3018.9167
$ws.Cells.Item(99, 12).Value = 2814.889
$ws.Cells.Item(99, 13).Value = -1520.9167
$ws.Cells.Item(99, 14).Value = -5810.889

$ws.Cells.Item(104, 8).Value = 49999.5
$ws.Cells.Item(104, 10).Value = 49999.5
$ws.Cells.Item(104, 12).Value = 49999.5
$ws.Cells.Item(104, 14).Value = -55241.5

$ws.Cells.Item(105, 8).Value = 671.5714
$ws.Cells.Item(105, 9).Value = 616.8333
$ws.Cells.Item(105, 11).Value = 616.8333
$ws.Cells.Item(105, 13).Value = 1130.1667

$ws.Cells.Item(106, 8).Value = 62194.332
$ws.Cells.Item(106, 10).Value = 62194.332
$ws.Cells.Item(106, 12).Value = 62194.332
$ws.Cells.Item(106, 14).Value = -64718.332

$ws.Cells.Item(122, 8).Value = 2924.6538
$ws.Cells.Item(122, 9).Value = 3291.8948
$ws.Cells.Item(122, 11).Value = 9875.6844
$ws.Cells.Item(122, 13).Value = -7425.6844

$ws.Cells.Item(126, 8).Value = 2931.476
$ws.Cells.Item(126, 9).Value = 3018.9167
$ws.Cells.Item(126, 10).Value = 2814.889
$ws.Cells.Item(126, 11).Value = 9056.750100000001
$ws.Cells.Item(126, 12).Value = 8444.667000000001
$ws.Cells.Item(126, 13).Value = -6586.750100000001
$ws.Cells.Item(126, 14).Value = -13384.667

$ws.Cells.Item(132, 8).Value = 3381
$ws.Cells.Item(132, 9).Value = 1763.25
$ws.Cells.Item(132, 11).Value = 5289.75
$ws.Cells.Item(132, 13).Value = -2759.75

$ws.Cells.Item(134, 8).Value = 8502
$ws.Cells.Item(134, 9).Value = 10246.5
$ws.Cells.Item(134, 10).Value = 5013
$ws.Cells.Item(134, 11).Value = 30739.5
$ws.Cells.Item(134, 12).Value = 15039
$ws.Cells.Item(134, 13).Value = -28204.5
$ws.Cells.Item(134, 14).Value = -20109

$ws.Cells.Item(136, 8).Value = 2026.15
$ws.Cells.Item(136, 9).Value = 1166.6666
$ws.Cells.Item(136, 10).Value = 2729.3635
$ws.Cells.Item(136, 11).Value = 3499.9998
$ws.Cells.Item(136, 12).Value = 8188.0905
$ws.Cells.Item(136, 13).Value = -949.9998000000001
$ws.Cells.Item(136, 14).Value = -13288.0905

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(42, 8).Value = 12900
$ws.Cells.Item(42, 9).Value = 11500
$ws.Cells.Item(42, 10).Value = 15000
$ws.Cells.Item(42, 11).Value = 34500
$ws.Cells.Item(42, 12).Value = 45000
$ws.Cells.Item(42, 13).Value = -33966
$ws.Cells.Item(42, 14).Value = -46068

$ws.Cells.Item(44, 8).Value = 0
$ws.Cells.Item(44, 9).Value = 0
$ws.Cells.Item(44, 10).Value = 0
$ws.Cells.Item(44, 11).Value = 0
$ws.Cells.Item(44, 12).Value = 0
$ws.Cells.Item(44, 13).ClearContents()
$ws.Cells.Item(44, 14).ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 3713.6667
$ws.Cells.Item(102, 9).Value = 3689.15
$ws.Cells.Item(102, 10).Value = 3762.7
$ws.Cells.Item(102, 11).Value = 3689.15
$ws.Cells.Item(102, 12).Value = 3762.7
$ws.Cells.Item(102, 13).Value = -2067.15
$ws.Cells.Item(102, 14).Value = -7006.7

$ws.Cells.Item(132, 8).Value = 6727.6113
$ws.Cells.Item(132, 9).Value = 6591.7144
$ws.Cells.Item(132, 10).Value = 7203.25
$ws.Cells.Item(132, 11).Value = 19775.1432
$ws.Cells.Item(132, 12).Value = 21609.75
$ws.Cells.Item(132, 13).Value = -17245.1432
$ws.Cells.Item(132, 14).Value = -26669.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 5429.5713
$ws.Cells.Item(7, 9).Value = 2752
$ws.Cells.Item(7, 10).Value = 8999.666999999999
$ws.Cells.Item(7, 11).Value = 2752
$ws.Cells.Item(7, 12).Value = 8999.666999999999
$ws.Cells.Item(7, 13).Value = -2640
$ws.Cells.Item(7, 14).Value = -9223.666999999999

$ws.Cells.Item(22, 8).Value = 6994035
$ws.Cells.Item(22, 9).Value = 12987596
$ws.Cells.Item(22, 10).Value = 1548
$ws.Cells.Item(22, 11).Value = 12987596
$ws.Cells.Item(22, 12).Value = 1548
$ws.Cells.Item(22, 13).Value = -12987301
$ws.Cells.Item(22, 14).Value = -2138

$ws.Cells.Item(27, 8).Value = 6994035
$ws.Cells.Item(27, 9).Value = 12987596
$ws.Cells.Item(27, 10).Value = 1548
$ws.Cells.Item(27, 11).Value = 12987596
$ws.Cells.Item(27, 12).Value = 1548
$ws.Cells.Item(27, 13).Value = -12987489
$ws.Cells.Item(27, 14).Value = -1762

$ws.Cells.Item(126, 8).Value = 5429.5713
$ws.Cells.Item(126, 9).Value = 2752
$ws.Cells.Item(126, 10).Value = 8999.666999999999
$ws.Cells.Item(126, 11).Value = 8256
$ws.Cells.Item(126, 12).Value = 26999.001
$ws.Cells.Item(126, 13).Value = -5786
$ws.Cells.Item(126, 14).Value = -31939.001

$ws.Cells.Item(132, 8).Value = 3007.7334
$ws.Cells.Item(132, 9).Value = 2473.7273
$ws.Cells.Item(132, 11).Value = 7421.1819
$ws.Cells.Item(132, 13).Value = -4891.1819

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(45, 8).Value = 12391.777
$ws.Cells.Item(45, 9).Value = 4468
$ws.Cells.Item(45, 10).Value = 14655.714
$ws.Cells.Item(45, 11).Value = 4468
$ws.Cells.Item(45, 12).Value = 14655.714
$ws.Cells.Item(45, 13).Value = -3977
$ws.Cells.Item(45, 14).Value = -15637.714

$ws.Cells.Item(81, 8).Value = 2015.2142
$ws.Cells.Item(81, 9).Value = 1383
$ws.Cells.Item(81, 10).Value = 4333.3335
$ws.Cells.Item(81, 11).Value = 2766
$ws.Cells.Item(81, 12).Value = 8666.666999999999
$ws.Cells.Item(81, 13).Value = -1705
$ws.Cells.Item(81, 14).Value = -10788.667

$ws.Cells.Item(84, 8).Value = 2015.2142
$ws.Cells.Item(84, 9).Value = 1383
$ws.Cells.Item(84, 10).Value = 4333.3335
$ws.Cells.Item(84, 11).Value = 13830
$ws.Cells.Item(84, 12).Value = 43333.335
$ws.Cells.Item(84, 13).Value = -8526
$ws.Cells.Item(84, 14).Value = -53941.335

$ws.Cells.Item(107, 8).Value = 1034.6
$ws.Cells.Item(107, 10).Value = 1091.6666
$ws.Cells.Item(107, 12).Value = 3274.9998
$ws.Cells.Item(107, 14).Value = -7114.9998

$ws.Cells.Item(126, 8).Value = 6030.1763
$ws.Cells.Item(126, 9).Value = 6358.143
$ws.Cells.Item(126, 11).Value = 19074.429
$ws.Cells.Item(126, 13).Value = -16604.429
